$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 189: 30-09-2021
$ws.Cells.Item(189, 1).Formula = '="30-09-2021"'
$ws.Cells.Item(189, 1).Copy()
$ws.Cells.Item(189, 1).PasteSpecial(-4163)
$ws.Cells.Item(189, 2).Value = 0.16
$ws.Cells.Item(189, 3).Value = 0.19
$ws.Cells.Item(189, 4).Value = 0.01
$ws.Cells.Item(189, 5).Value = 0.23
$ws.Cells.Item(189, 6).Value = 0.32

# Row 190: 01-10-2021
$ws.Cells.Item(190, 1).Formula = '="01-10-2021"'
$ws.Cells.Item(190, 1).Copy()
$ws.Cells.Item(190, 1).PasteSpecial(-4163)
$ws.Cells.Item(190, 2).Value = 0.15
$ws.Cells.Item(190, 3).Value = 0.24
$ws.Cells.Item(190, 4).Value = 0.01
$ws.Cells.Item(190, 5).Value = 0.08
$ws.Cells.Item(190, 6).Value = 0.35
